$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = '@'
    $range.Value = $value
    $range.Style = 'Normal'
}

# Row 2
Set-TextValue $ws.Range("D2") '43.006.73'
Set-TextValue $ws.Range("E2") '  -7.03%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.549.88'
Set-TextValue $ws.Range("E3") '  -2.68%  '

# Row 4
Set-TextValue $ws.Range("E4") '  +0.03%  '

# Row 5
Set-TextValue $ws.Range("D5") '298.39'
Set-TextValue $ws.Range("E5") '  -4.10%  '

# Row 6
Set-TextValue $ws.Range("D6") '92.37'
Set-TextValue $ws.Range("E6") '  -7.28%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.574'
Set-TextValue $ws.Range("E7") '  -3.91%  '

# Row 8
Set-TextValue $ws.Range("E8") '  +0.00%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.551'
Set-TextValue $ws.Range("E9") '  -5.68%  '

# Row 10
Set-TextValue $ws.Range("D10") '36.10'
Set-TextValue $ws.Range("E10") '  -7.63%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.0809'
Set-TextValue $ws.Range("E11") '  -4.38%  '

# Row 12
Set-TextValue $ws.Range("D12") '7.67'
Set-TextValue $ws.Range("E12") '  -5.54%  '

# Row 13
Set-TextValue $ws.Range("E13") '  +1.19%  '

# Row 14
Set-TextValue $ws.Range("D14") '2.939.77'
Set-TextValue $ws.Range("E14") '  -2.66%  '

# Row 15
Set-TextValue $ws.Range("D15") '2.552.85'
Set-TextValue $ws.Range("E15") '  -2.23%  '

# Row 16
Set-TextValue $ws.Range("E16") '  -5.31%  '

# Row 17
Set-TextValue $ws.Range("D17") '14.18'
Set-TextValue $ws.Range("E17") '  -5.00%  '

# Row 18
Set-TextValue $ws.Range("D18") '43.060.69'
Set-TextValue $ws.Range("E18") '  -7.41%  '

# Row 19
Set-TextValue $ws.Range("D19") '12.82'
Set-TextValue $ws.Range("E19") '  -0.30%  '

# Row 20
Set-TextValue $ws.Range("D20") '0.0₃0979'
Set-TextValue $ws.Range("E20") '  -4.19%  '

# Row 21
Set-TextValue $ws.Range("D21") '6.65'
Set-TextValue $ws.Range("E21") '  -2.30%  '

# Row 22
Set-TextValue $ws.Range("D22") '71.75'
Set-TextValue $ws.Range("E22") '  -2.36%  '

# Row 23
Set-TextValue $ws.Range("D23") '260.38'
Set-TextValue $ws.Range("E23") '  -11.30%  '

# Row 24
Set-TextValue $ws.Range("D24") '2.91'
Set-TextValue $ws.Range("E24") '  -5.16%  '

# Row 25
Set-TextValue $ws.Range("B25") 'EthereumClassic'
Set-TextValue $ws.Range("C25") 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D25") '29.56'
Set-TextValue $ws.Range("E25") '  -0.71%  '

# Row 26
Set-TextValue $ws.Range("B26") 'ImmutableX'
Set-TextValue $ws.Range("C26") 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D26") '2.14'
Set-TextValue $ws.Range("E26") '  -4.77%  '

# Row 27
Set-TextValue $ws.Range("E27") '  +0.07%  '

# Row 28
Set-TextValue $ws.Range("D28") '10.05'
Set-TextValue $ws.Range("E28") '  -7.46%  '

# Row 29
Set-TextValue $ws.Range("B29") 'Toncoin'
Set-TextValue $ws.Range("C29") 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D29") '2.13'
Set-TextValue $ws.Range("E29") '  -4.08%  '

# Row 30
Set-TextValue $ws.Range("B30") 'InjectiveProtocol'
Set-TextValue $ws.Range("C30") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D30") '36.76'
Set-TextValue $ws.Range("E30") '  -6.48%  '

# Row 31
Set-TextValue $ws.Range("E31") '  -5.17%  '

# Row 32
Set-TextValue $ws.Range("D32") '154.61'
Set-TextValue $ws.Range("E32") '  -2.48%  '

# Row 33
Set-TextValue $ws.Range("E33") '  -3.15%  '

# Row 34
Set-TextValue $ws.Range("D34") '3.39'
Set-TextValue $ws.Range("E34") '  -6.08%  '

# Row 35
Set-TextValue $ws.Range("E35") '  -2.59%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.0798'
Set-TextValue $ws.Range("E36") '  -5.42%  '

# Row 37
Set-TextValue $ws.Range("E37") '  -5.87%  '

# Row 38
Set-TextValue $ws.Range("E38") '  -3.36%  '

# Row 39
Set-TextValue $ws.Range("D39") '23.33'
Set-TextValue $ws.Range("E39") '  +7.86%  '

# Row 40
Set-TextValue $ws.Range("D40") '16.50'
Set-TextValue $ws.Range("E40") '  +4.17%  '

# Row 41
Set-TextValue $ws.Range("D41") '3.46'
Set-TextValue $ws.Range("E41") '  -3.68%  '

# Row 42
Set-TextValue $ws.Range("E42") '  -5.52%  '

# Row 43
Set-TextValue $ws.Range("E43") '  -4.23%  '

# Row 44
Set-TextValue $ws.Range("D44") '2.069.24'
Set-TextValue $ws.Range("E44") '  -2.71%  '

# Row 45
Set-TextValue $ws.Range("D45") '0.999'
Set-TextValue $ws.Range("E45") '  +0.02%  '

# Row 46
Set-TextValue $ws.Range("D46") '85.68'
Set-TextValue $ws.Range("E46") '  -11.94%  '

# Row 47
Set-TextValue $ws.Range("E47") '  +2.87%  '

# Row 48
Set-TextValue $ws.Range("D48") '2.795.79'
Set-TextValue $ws.Range("E48") '  -2.70%  '

# Row 49
Set-TextValue $ws.Range("B49") 'FraxShare'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D49") '8.76'
Set-TextValue $ws.Range("E49") '  -8.19%  '

# Row 50
Set-TextValue $ws.Range("B50") 'Stacks'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D50") '1.70'
Set-TextValue $ws.Range("E50") '  -2.91%  '

# Row 51
Set-TextValue $ws.Range("D51") '104.25'
Set-TextValue $ws.Range("E51") '  -5.28%  '
